# Apply the "config reader" update described by the commit:
#  - Settings sheet: reshuffle the email/share-drive rows, add a new
#    "emailName" row and a new "fromEmailAddress" row, drop the old
#    "_orig" hyperlinked rows/hyperlinks.
#  - Constants sheet: add a new MediumTimeSpan constant row.
#  - View state: Settings becomes the active tab/sheet; Constants keeps
#    its own selection but is no longer the active tab.
#  - Constants sheet gets an explicit portrait page setup.

$wb = $excel.ActiveWorkbook
$settings  = $wb.Worksheets.Item("Settings")
$constants = $wb.Worksheets.Item("Constants")

# ---------------------------------------------------------------------
# Settings sheet (sheet1)
# ---------------------------------------------------------------------

# Remove the two existing hyperlinks (they live on B11/B12 before the
# reshuffle below moves everything around).
$settings.Range("B11").Hyperlinks.Delete()
$settings.Range("B12").Hyperlinks.Delete()

# Row 10: was shareDrive_orig -> becomes the new "emailName" row.
$settings.Range("A10").Value = "emailName"
$settings.Range("B10").Value = "Truman  GSA"
$settings.Range("C10").Value = '"Name" to use for return email'

# Row 11: was errorEmailAddress_orig -> now blank.
$settings.Range("A11:C11").Clear()

# Row 12: shareDrive (moved down from row 10). B12 previously held the
# hyperlink-styled successEmailAddress_orig value, so reset its style.
$settings.Range("A12").Value = "shareDrive"
$settings.Range("B12").Style = "Normal"
$settings.Range("B12").Value = "\\E04BMV-CIFS02.ent.ds.gsa.gov\R02_FSS_Shares$\Files\Shares\TrumanRPA\"
$settings.Range("C12").Value = "Shared drive to be used"

# Row 13: errorEmailAddress (moved down from row 11), keeps the
# hyperlink-like style but no longer has a live hyperlink.
$settings.Range("A13").Value = "errorEmailAddress"
$settings.Range("B13").Value = "truman.00corp.testers@gsa.gov"
$settings.Range("B13").Style = "Hyperlink"
$settings.Range("C13").Value = "email address to send error email"

# Row 14: successEmailAddress (moved down from row 12).
$settings.Range("A14").Value = "successEmailAddress"
$settings.Range("B14").Value = "truman.00corp.testers@gsa.gov"
$settings.Range("B14").Style = "Hyperlink"
$settings.Range("C14").Value = "email address to send sucess email"

# Row 15: brand new fromEmailAddress row.
$settings.Range("A15").Value = "fromEmailAddress"
$settings.Range("B15").Value = "truman.00corp.testers@gsa.gov"
$settings.Range("B15").Style = "Hyperlink"
$settings.Range("C15").Value = 'email address to use in the "from" field'

# Extend the used range down to row 992 (two more blank formatted rows).
$settings.Rows.Item(991).RowHeight = 14.25
$settings.Rows.Item(992).RowHeight = 14.25

# ---------------------------------------------------------------------
# Constants sheet (sheet2)
# ---------------------------------------------------------------------

# New row 29: MediumTimeSpan constant.
$constants.Range("A29").Value = "MediumTimeSpan"
$constants.Range("B29").NumberFormat = "@"
$constants.Range("B29").Value = "00:00:02"
$constants.Range("C29").Value = "2 second time span"

# Give the Constants sheet an explicit portrait page setup.
$constants.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# View / selection state
#   Settings becomes the active sheet/tab with A15 selected.
#   Constants keeps C29 selected but is no longer the active tab.
# ---------------------------------------------------------------------

$constants.Activate()
$constants.Range("C29").Select()

$settings.Activate()
$settings.Range("A15").Select()
